# Plumbing Client v1 update
# Replaces the "AVS International school" Q&A transcript rows with the new
# "Pete's Plumbing" transcript rows, restyles the header row, clears the
# now-unused rows, and updates row heights / sheet selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rsquo = [char]0x2019   # curly right single quote used throughout the new copy

# ---------------------------------------------------------------------
# 1. Clear the old data from rows 8-17 first (columns C/D in rows 2-17
#    were already empty). Clearing before writing the new rows lets the
#    shared-string table shrink naturally to just the strings still in
#    use once we're done (16 unique strings).
# ---------------------------------------------------------------------
$ws.Range("A8:D17").ClearContents()

# ---------------------------------------------------------------------
# 2. Write the new transcript content into rows 2-7 (columns A & B).
#    Column headers (row 1) and columns C/D for these rows keep their
#    existing (empty / header) values.
# ---------------------------------------------------------------------
$ws.Cells.Item(2,1).Value2 = "plumbing_intro.mp3"
$ws.Cells.Item(2,2).Value2 = "G" + $rsquo + "day! You" + $rsquo + "ve reached Pete's Plumbing. How can we help you today?`n"

$ws.Cells.Item(3,1).Value2 = "in_business_how_long.mp3"
$ws.Cells.Item(3,2).Value2 = "We" + $rsquo + "ve been doing this for 7 years now and got loads of happy clients and repeat work. You" + $rsquo + "ll be in safe hands.`n"

$ws.Cells.Item(4,1).Value2 = "services_offered.mp3"
$ws.Cells.Item(4,2).Value2 = "We handle blocked drains, leaking taps, toilet repairs, hot water issues, and a range of other plumbing issues. We also do gas fitting, pipe relining, and kitchen or bathroom plumbing. What" + $rsquo + "s the issue you" + $rsquo + "re facing right now?`n"

$ws.Cells.Item(5,1).Value2 = "available_hours.mp3"
$ws.Cells.Item(5,2).Value2 = "We" + $rsquo + "re available Monday through Saturday, 8am to 6pm. After-hours emergency support is also available for an extra service fee, so yeah, give us a ring anytime.`n"

$ws.Cells.Item(6,1).Value2 = "pricing.mp3"
$ws.Cells.Item(6,2).Value2 = "Our pricing usually starts at `$98 for standard service calls. We" + $rsquo + "ll provide a full quote after understanding the job better.`n"

$ws.Cells.Item(7,1).Value2 = "ask_time_day.mp3"
$ws.Cells.Item(7,2).Value2 = "We" + $rsquo + "ve got a few open slots this week. What time and day works for you?"

# ---------------------------------------------------------------------
# 3. Row heights - new rows are noticeably shorter than the old ones.
# ---------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 13.5
$ws.Rows.Item(2).RowHeight = 39
$ws.Rows.Item(3).RowHeight = 64.5
$ws.Rows.Item(4).RowHeight = 64.5
$ws.Rows.Item(5).RowHeight = 64.5
$ws.Rows.Item(6).RowHeight = 51.75
for ($r = 7; $r -le 17; $r++) {
  $ws.Rows.Item($r).RowHeight = 13.5
}

# ---------------------------------------------------------------------
# 4. Header row (row 1) gets a new bold white-on-green look.
# ---------------------------------------------------------------------
$ws.Range("A1:D1").Font.Bold = $true
$ws.Range("A1:D1").Font.ThemeColor = 2
$ws.Range("A1:D1").Font.Name = "Roboto"

# ---------------------------------------------------------------------
# 5. Sheet view - scroll back to the top and move the selection.
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 1
[void]$ws.Range("B11").Select()

Write-Output "Plumbing Client v1 update applied"
